# Fix the typo/grammar in the dashboard description on slide 7:
#   "will allows ... pendemic"  ->  "will allow ... pandemic"
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item("Text Placeholder 9")
$tr = $sh.TextFrame.TextRange

# --- Fix #1: "will allows " -> "will allow " ---------------------------
$full = $tr.Text
$idx  = $full.IndexOf("will allows")
if ($idx -ge 0) {
    $fix1 = $tr.Characters($idx + 1, 12)
    $fix1.Text = "will allow "
}

# --- Fix #2: "pendemic" -> "pandemic" (merge the rest of the sentence) -
$full = $tr.Text
$idx  = $full.IndexOf("us to better")
if ($idx -ge 0) {
    $tail = $tr.Characters($idx + 1, $full.Length - $idx)
    $tail.Text = "us to better understand labor trends caused by major events such as the Covid-19 pandemic and the 2008 recession. "
}
